$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.031.79"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.885.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7352"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3164"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07176"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("E10").Value = "  -1.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08319"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.392"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.886.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.149"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.045.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "248.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007841"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.143.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.890"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1566"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.274"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.046"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.474"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.576"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.31%  "

$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.195"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05318"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "

$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7687"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9980"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4577"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.036"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.086.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8739"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9996"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.560"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.537"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.039.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
